# Auto-generated edit script applying the Jenova_Profits market-data refresh
# diff across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting-profit tables.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 91250
$ws.Range("J87").Value = 91250
$ws.Range("L87").Value = 91250
$ws.Range("N87").Value = -93746
$ws.Range("H90").Value = 91250
$ws.Range("J90").Value = 91250
$ws.Range("L90").Value = 273750
$ws.Range("N90").Value = -286230
$ws.Range("H112").Value = 2999.8076
$ws.Range("J112").Value = 2999.8076
$ws.Range("L112").Value = 8999.4228
$ws.Range("N112").Value = -11215.4228
$ws.Range("H127").Value = 1714.75
$ws.Range("I127").Value = 686.6667
$ws.Range("K127").Value = 2060.0001
$ws.Range("M127").Value = 2899.9999
$ws.Range("H133").Value = 41779.223
$ws.Range("J133").Value = 41779.223
$ws.Range("L133").Value = 41779.223
$ws.Range("N133").Value = -51899.223
$ws.Range("H138").Value = 4785.0796
$ws.Range("J138").Value = 5757.821
$ws.Range("L138").Value = 17273.463
$ws.Range("N138").Value = -27553.463
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 3099.25
$ws.Range("I141").Value = 2709.7222
$ws.Range("K141").Value = 8129.1666
$ws.Range("M141").Value = -2949.1666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 9999
$ws.Range("I3").Value = 9999
$ws.Range("K3").Value = 9999
$ws.Range("M3").Value = -9884
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H32").Value = 3771.514
$ws.Range("I32").Value = 3436.1738
$ws.Range("K32").Value = 3436.1738
$ws.Range("M32").Value = -3149.1738
$ws.Range("H61").Value = 3008.76
$ws.Range("I61").Value = 2909.5217
$ws.Range("J61").Value = 4150
$ws.Range("K61").Value = 2909.5217
$ws.Range("L61").Value = 4150
$ws.Range("M61").Value = -2697.5217
$ws.Range("N61").Value = -4574
$ws.Range("H63").Value = 5653.6665
$ws.Range("I63").Value = 2549.1428
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 2549.1428
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -1863.1428
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 5653.6665
$ws.Range("I66").Value = 2549.1428
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 12745.714
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -9313.714
$ws.Range("N66").Value = -56864
$ws.Range("H74").Value = 1569.6538
$ws.Range("I74").Value = 1514.2174
$ws.Range("K74").Value = 1514.2174
$ws.Range("M74").Value = -640.2174
$ws.Range("H77").Value = 1569.6538
$ws.Range("I77").Value = 1514.2174
$ws.Range("K77").Value = 7571.087
$ws.Range("M77").Value = -3203.087
$ws.Range("H92").Value = 50000000
$ws.Range("I92").Value = 50000000
$ws.Range("K92").Value = 50000000
$ws.Range("M92").Value = -49997504
$ws.Range("H125").Value = 76189.664
$ws.Range("J125").Value = 76189.664
$ws.Range("L125").Value = 76189.664
$ws.Range("N125").Value = -86029.664
$ws.Range("H136").Value = 3008.76
$ws.Range("I136").Value = 2909.5217
$ws.Range("J136").Value = 4150
$ws.Range("K136").Value = 8728.5651
$ws.Range("L136").Value = 12450
$ws.Range("M136").Value = -6178.5651
$ws.Range("N136").Value = -17550
$ws.Range("H138").Value = 75000
$ws.Range("J138").Value = 75000
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 112268.664
$ws.Range("I105").Value = 112268.664
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 112268.664
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -110521.664
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 456728.97
$ws.Range("I107").Value = 1676.375
$ws.Range("J107").Value = 1670202.5
$ws.Range("K107").Value = 1676.375
$ws.Range("L107").Value = 1670202.5
$ws.Range("M107").Value = 243.625
$ws.Range("N107").Value = -1674042.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38980.895
$ws.Range("I31").Value = 1109.75
$ws.Range("K31").Value = 1109.75
$ws.Range("M31").Value = -814.75
$ws.Range("H34").Value = 38980.895
$ws.Range("I34").Value = 1109.75
$ws.Range("K34").Value = 1109.75
$ws.Range("M34").Value = -907.75
$ws.Range("H58").Value = 4556
$ws.Range("I58").Value = 5344.4
$ws.Range("J58").Value = 614
$ws.Range("K58").Value = 5344.4
$ws.Range("L58").Value = 614
$ws.Range("M58").Value = -5141.4
$ws.Range("N58").Value = -1020
$ws.Range("H99").Value = 4312.647
$ws.Range("I99").Value = 2290.125
$ws.Range("J99").Value = 6110.4443
$ws.Range("K99").Value = 2290.125
$ws.Range("L99").Value = 6110.4443
$ws.Range("M99").Value = -792.125
$ws.Range("N99").Value = -9106.444299999999
$ws.Range("H122").Value = 4163.2354
$ws.Range("I122").Value = 3661.625
$ws.Range("J122").Value = 4609.1113
$ws.Range("K122").Value = 10984.875
$ws.Range("L122").Value = 13827.3339
$ws.Range("M122").Value = -8534.875
$ws.Range("N122").Value = -18727.3339
$ws.Range("H126").Value = 4312.647
$ws.Range("I126").Value = 2290.125
$ws.Range("J126").Value = 6110.4443
$ws.Range("K126").Value = 6870.375
$ws.Range("L126").Value = 18331.3329
$ws.Range("M126").Value = -4400.375
$ws.Range("N126").Value = -23271.3329
$ws.Range("H136").Value = 4556
$ws.Range("I136").Value = 5344.4
$ws.Range("J136").Value = 614
$ws.Range("K136").Value = 16033.2
$ws.Range("L136").Value = 1842
$ws.Range("M136").Value = -13483.2
$ws.Range("N136").Value = -6942

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 11739692
$ws.Range("J11").Value = 950.9
$ws.Range("L11").Value = 2852.7
$ws.Range("N11").Value = -3132.7
$ws.Range("H63").Value = 1900.3846
$ws.Range("I63").Value = 1610.4546
$ws.Range("J63").Value = 3495
$ws.Range("K63").Value = 4831.3638
$ws.Range("L63").Value = 10485
$ws.Range("M63").Value = -4082.3638
$ws.Range("N63").Value = -11983
$ws.Range("H64").Value = 166667220
$ws.Range("I64").Value = 166667220
$ws.Range("K64").Value = 500001660
$ws.Range("M64").Value = -500001390
$ws.Range("H66").Value = 1900.3846
$ws.Range("I66").Value = 1610.4546
$ws.Range("J66").Value = 3495
$ws.Range("K66").Value = 14494.0914
$ws.Range("L66").Value = 31455
$ws.Range("M66").Value = -10750.0914
$ws.Range("N66").Value = -38943
$ws.Range("H67").Value = 166667220
$ws.Range("I67").Value = 166667220
$ws.Range("K67").Value = 500001660
$ws.Range("M67").Value = -500000724
$ws.Range("H70").Value = 127169.875
$ws.Range("I70").Value = 127169.875
$ws.Range("K70").Value = 381509.625
$ws.Range("M70").Value = -381194.625
$ws.Range("H73").Value = 127169.875
$ws.Range("I73").Value = 127169.875
$ws.Range("K73").Value = 381509.625
$ws.Range("M73").Value = -380417.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 10000
$ws.Range("K5").Value = 10000
$ws.Range("M5").Value = -9888
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 4787.4
$ws.Range("I122").Value = 3644.8572
$ws.Range("K122").Value = 10934.5716
$ws.Range("M122").Value = -8484.571599999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1889.8462
$ws.Range("I46").Value = 2296.2222
$ws.Range("J46").Value = 975.5
$ws.Range("K46").Value = 2296.2222
$ws.Range("L46").Value = 975.5
$ws.Range("M46").Value = -2108.2222
$ws.Range("N46").Value = -1351.5
$ws.Range("H68").Value = 168416.83
$ws.Range("I68").Value = 867
$ws.Range("J68").Value = 335966.66
$ws.Range("K68").Value = 867
$ws.Range("L68").Value = 335966.66
$ws.Range("M68").Value = -118
$ws.Range("N68").Value = -337464.66
$ws.Range("H71").Value = 168416.83
$ws.Range("I71").Value = 867
$ws.Range("J71").Value = 335966.66
$ws.Range("K71").Value = 4335
$ws.Range("L71").Value = 1679833.3
$ws.Range("M71").Value = -591
$ws.Range("N71").Value = -1687321.3
$ws.Range("H100").Value = 3098
$ws.Range("I100").Value = 2800
$ws.Range("K100").Value = 2800
$ws.Range("M100").Value = -2259

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 38000
$ws.Range("J54").Value = 38000
$ws.Range("L54").Value = 38000
$ws.Range("N54").Value = -39040
$ws.Range("H81").Value = 2919.5652
$ws.Range("I81").Value = 1492.4
$ws.Range("J81").Value = 12434
$ws.Range("K81").Value = 2984.8
$ws.Range("L81").Value = 24868
$ws.Range("M81").Value = -1923.8
$ws.Range("N81").Value = -26990
$ws.Range("H84").Value = 2919.5652
$ws.Range("I84").Value = 1492.4
$ws.Range("J84").Value = 12434
$ws.Range("K84").Value = 14924
$ws.Range("L84").Value = 124340
$ws.Range("M84").Value = -9620
$ws.Range("N84").Value = -134948
$ws.Range("H132").Value = 44693.96
$ws.Range("I132").Value = 3646.5
$ws.Range("K132").Value = 10939.5
$ws.Range("M132").Value = -8409.5

Write-Host "Applied market-data refresh: 246 cells set, 4 cells cleared"